# Update "想去人数" (want-to-go count) figures across sheets to match the
# newly generated data snapshot (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 807
$wsExpo.Range("F4").Value = 1096
$wsExpo.Range("F7").Value = 215
$wsExpo.Range("F8").Value = 390
$wsExpo.Range("F11").Value = 506
$wsExpo.Range("F12").Value = 535
$wsExpo.Range("F14").Value = 12603
$wsExpo.Range("F15").Value = 5166
$wsExpo.Range("F16").Value = 5516

# --- Sheet "演出" (Performance) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F3").Value = 28

# --- Sheet "全部类型" (All Types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 807
$wsAll.Range("F5").Value = 1097
$wsAll.Range("F8").Value = 215
$wsAll.Range("F9").Value = 390
$wsAll.Range("F12").Value = 506
$wsAll.Range("F13").Value = 535
$wsAll.Range("F15").Value = 12603
$wsAll.Range("F16").Value = 28
$wsAll.Range("F18").Value = 5166
$wsAll.Range("F19").Value = 5516
